$d = $word.ActiveDocument

$replacements = @(
    @("732×5=", "961×2="),
    @("932×6=", "744×9="),
    @("583×7=", "785×7="),
    @("520×7=", "891×8="),
    @("433×7=", "313×7="),
    @("401×4=", "658×4="),
    @("894×4=", "632×4="),
    @("933×5=", "881×5="),
    @("200×5=", "764×8="),
    @("196×4=", "756×4="),
    @("427×2=", "665×3="),
    @("225×3=", "481×6="),
    @("929×4=", "607×7="),
    @("953×2=", "144×7="),
    @("391×5=", "552×5="),
    @("583×8=", "751×3="),
    @("980×6=", "918×6="),
    @("780×4=", "181×7="),
    @("807×9=", "680×9="),
    @("401×9=", "421×7="),
    @("806×5=", "720×2="),
    @("908×8=", "431×6="),
    @("158×8=", "422×8="),
    @("279×5=", "684×9="),
    @("350×6=", "387×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
